$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Llama 3.2 1B" row (row 4) and "Llama 3.2m 3B" row (originally row 10,
# now row 9 after the first deletion).
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(9).Delete()

# Update the selection to match the target view state.
$ws.Range("E21:E22").Select()
